# Restore the presentation's (Integral / "Red Violet") theme colour scheme
# back to the stock Office theme colours ("Office" colour scheme: dk1,
# lt1, dk2, lt2, accent1-6, hlink, folHlink), matching the colours already
# present in the deck's notes-master theme part.
#
# PowerPoint COM represents a theme's twelve colour slots through
# ThemeColorScheme.Colors(1..12).RGB (index order: dk1, lt1, dk2, lt2,
# accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink).
# .RGB takes/returns a single OLE colour value packed as 0x00BBGGRR, so a
# small helper converts the familiar RRGGBB hex strings used in the theme
# XML into that packed form.

function Convert-HexToOleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target values: the stock "Office" colour scheme (dk1..folHlink), in
# ThemeColorScheme.Colors() slot order.
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = Convert-HexToOleColor $officeColors[$i - 1]
}
